$d = $word.ActiveDocument

# Locate the paragraph that contains the original sentence.
$oldText = "Check to see if the other values are reversed versions of the current one."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "$oldText*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    $r = $target.Range

    # Pull the paragraph's own opening <w:p ...> tag (with its original
    # paraId/rsid attributes) out of its WordOpenXML so the replacement
    # keeps them intact.
    $owx = $r.WordOpenXML
    $openTag = "<w:p>"
    if ($owx -match '(<w:p[ />][^>]*>)') {
        $openTag = $matches[1]
    }

    $run1 = "<w:r><w:t>Create an empty set. See if the value of the words list has a reverse in the set</w:t></w:r>"
    $run2 = "<w:r><w:t xml:space=`"preserve`">. </w:t></w:r>"
    $run3 = "<w:r><w:t>If there is, print the match out. Then add the list value to the set so it can be paired later (provided there are pairs).</w:t></w:r>"
    $run4 = "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"

    $newParaBody = "$openTag$run1$run2$run3$run4</w:p>"

    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $newParaBody + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml)
}
